$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values (Strike# -> K) for rows 2-19, replacing the previous
# Strike# counts with the recalculated K values.
$kValues = @{
    2  = 8
    3  = 6
    4  = 4
    5  = 1
    6  = 4
    7  = 4
    8  = 5
    9  = 9
    10 = 2
    11 = 1
    12 = 0
    13 = 7
    14 = 6
    15 = 5
    16 = 6
    17 = 4
    18 = 4
    19 = 5
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
